$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cryptocurrency market data (coin name/link swaps, prices, and 1h volume %)

$ws.Cells.Item(2, 4).Value = "35.376.05"
$ws.Cells.Item(2, 5).Value = "  -3.98%  "

$ws.Cells.Item(3, 4).Value = "1.981.21"
$ws.Cells.Item(3, 5).Value = "  -5.45%  "

$ws.Cells.Item(4, 5).Value = "  +0.29%  "

$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "240.04"
$ws.Cells.Item(5, 5).Value = "  -2.06%  "

$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "0.632"
$ws.Cells.Item(6, 5).Value = "  -6.43%  "

$ws.Cells.Item(7, 2).Value = "USDC"
$ws.Cells.Item(7, 3).Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = "1.00"
$ws.Cells.Item(7, 5).Value = "  +0.11%  "

$ws.Cells.Item(8, 2).Value = "Solana"
$ws.Cells.Item(8, 3).Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = "56.07"
$ws.Cells.Item(8, 5).Value = "  +3.47%  "

$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = "59.24"
$ws.Cells.Item(9, 5).Value = "  -0.05%  "

$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = "0.356"
$ws.Cells.Item(10, 5).Value = "  -3.39%  "

$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = "0.0723"
$ws.Cells.Item(11, 5).Value = "  -6.01%  "

$ws.Cells.Item(12, 5).Value = "  -6.44%  "

$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = "0.892"
$ws.Cells.Item(13, 5).Value = "  -3.96%  "

$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = "14.21"
$ws.Cells.Item(14, 5).Value = "  -5.07%  "

$ws.Cells.Item(15, 4).Value = "2.277.40"
$ws.Cells.Item(15, 5).Value = "  -5.03%  "

$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = "5.24"
$ws.Cells.Item(16, 5).Value = "  -4.60%  "

$ws.Cells.Item(17, 4).Value = "1.979.97"
$ws.Cells.Item(17, 5).Value = "  -5.31%  "

$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = "17.11"
$ws.Cells.Item(18, 5).Value = "  -0.63%  "

$ws.Cells.Item(19, 4).Value = "35.293.33"
$ws.Cells.Item(19, 5).Value = "  -4.12%  "

$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = "69.80"
$ws.Cells.Item(20, 5).Value = "  -4.09%  "

$ws.Cells.Item(21, 4).Value = "0.0₃0831"
$ws.Cells.Item(21, 5).Value = "  -5.90%  "

$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = "230.84"
$ws.Cells.Item(22, 5).Value = "  -3.75%  "

$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = "4.99"
$ws.Cells.Item(23, 5).Value = "  -8.66%  "

$ws.Cells.Item(24, 5).Value = "  -0.13%  "

$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = "2.27"
$ws.Cells.Item(25, 5).Value = "  -5.36%  "

$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = "2.24"
$ws.Cells.Item(26, 5).Value = "  +3.86%  "

$ws.Cells.Item(27, 2).Value = "Monero"
$ws.Cells.Item(27, 3).Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = "162.76"
$ws.Cells.Item(27, 5).Value = "  -2.59%  "

$ws.Cells.Item(28, 2).Value = "Cosmos"
$ws.Cells.Item(28, 3).Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = "9.07"
$ws.Cells.Item(28, 5).Value = "  -6.02%  "

$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = "19.32"
$ws.Cells.Item(29, 5).Value = "  -8.32%  "

$ws.Cells.Item(30, 5).Value = "  -5.87%  "

$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = "1.13"
$ws.Cells.Item(31, 5).Value = "  -3.42%  "

$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = "4.74"
$ws.Cells.Item(32, 5).Value = "  -9.54%  "

$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = "0.0582"
$ws.Cells.Item(33, 5).Value = "  -4.40%  "

$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = "0.0895"
$ws.Cells.Item(34, 5).Value = "  +8.58%  "

$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = "4.23"
$ws.Cells.Item(35, 5).Value = "  -10.69%  "

$ws.Cells.Item(36, 5).Value = "  +0.25%  "

$ws.Cells.Item(37, 5).Value = "  -8.25%  "

$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = "1.81"
$ws.Cells.Item(38, 5).Value = "  -2.03%  "

$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = "4.83"
$ws.Cells.Item(39, 5).Value = "  -1.71%  "

$ws.Cells.Item(40, 5).Value = "  -7.70%  "

$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = "2.78"
$ws.Cells.Item(41, 5).Value = "  -2.88%  "

$ws.Cells.Item(42, 5).Value = "  -6.26%  "

$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = "1.07"
$ws.Cells.Item(43, 5).Value = "  -7.76%  "

$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = "0.0879"
$ws.Cells.Item(44, 5).Value = "  -8.65%  "

$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = "89.64"
$ws.Cells.Item(45, 5).Value = "  -7.27%  "

$ws.Cells.Item(46, 4).Value = "1.357.37"
$ws.Cells.Item(46, 5).Value = "  -3.90%  "

$ws.Cells.Item(47, 5).Value = "  -5.28%  "

$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = "15.38"
$ws.Cells.Item(48, 5).Value = "  -4.39%  "

$ws.Cells.Item(49, 5).Value = "  -0.98%  "

$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = "2.24"
$ws.Cells.Item(50, 5).Value = "  -8.25%  "

$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = "45.15"
$ws.Cells.Item(51, 5).Value = "  -1.94%  "
